# Fruta / hortaliza, semanal
# Insert two new weekly price records at the top of the Limón data block
# (rows 252-253), pushing the existing records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("252:253").Insert()

# New row 252: 1a amarillo
$ws.Range("A252").Value2 = 11
$ws.Range("B252").Value2 = "Vega Monumental Concepción"
$ws.Range("C252").Value2 = "Bíobío"
$ws.Range("D252").Value2 = 44511
$ws.Range("E252").Value2 = 8
$ws.Range("F252").Value2 = "Fruta"
$ws.Range("G252").Value2 = 100102
$ws.Range("H252").Value2 = "Cítricos"
$ws.Range("I252").Value2 = 100102003
$ws.Range("J252").Value2 = "Limón"
$ws.Range("K252").Value2 = "Sin especificar"
$ws.Range("L252").Value2 = "1a amarillo"
$ws.Range("M252").Value2 = 550
$ws.Range("N252").Value2 = 5500
$ws.Range("O252").Value2 = 6000
$ws.Range("P252").Value2 = 5727
$ws.Range("Q252").Value2 = "$/malla 16 kilos"
$ws.Range("R252").Value2 = "Región de O'Higgins"
$ws.Range("S252").Value2 = 358
$ws.Range("T252").Value2 = 16

# New row 253: 1a plateado
$ws.Range("A253").Value2 = 11
$ws.Range("B253").Value2 = "Vega Monumental Concepción"
$ws.Range("C253").Value2 = "Bíobío"
$ws.Range("D253").Value2 = 44511
$ws.Range("E253").Value2 = 8
$ws.Range("F253").Value2 = "Fruta"
$ws.Range("G253").Value2 = 100102
$ws.Range("H253").Value2 = "Cítricos"
$ws.Range("I253").Value2 = 100102003
$ws.Range("J253").Value2 = "Limón"
$ws.Range("K253").Value2 = "Sin especificar"
$ws.Range("L253").Value2 = "1a plateado"
$ws.Range("M253").Value2 = 250
$ws.Range("N253").Value2 = 7000
$ws.Range("O253").Value2 = 7500
$ws.Range("P253").Value2 = 7200
$ws.Range("Q253").Value2 = "$/malla 16 kilos"
$ws.Range("R253").Value2 = "Región de O'Higgins"
$ws.Range("S253").Value2 = 450
$ws.Range("T253").Value2 = 16

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
